$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the four long description cells with a single-paragraph version of the
# same text (the embedded line breaks were removed) and drop their wrap-text
# formatting back to normal.
$ws.Range("E3").Value = '<b>A CO JEŚLI NIE CHCĘ JESZCZE MIESZKAĆ W MIEŚCIE?</b><br><br><p>,,Ujemne saldo migracji na wieś odnotowano w grupach wieku 20–24 i 25–29 lat. Zjawisko to wiąże się z odpływem młodych ludzi do miast, by kontynuować edukację lub podjąć pracę. (...) Zmiana miejsca zamieszkania następująca po ukończeniu nauki w szkole średniej wiąże się z nauką na studiach i następnie pracą w mieście, do którego młodzi dorośli migrują. Ta nierejestrowana migracja znajduje odzwierciedlenie w zmianie zameldowania najczęściej dopiero w momencie wymuszenia tego zdarzenia na przykład zakupem mieszkania.”</p><p align=”right”>,,Sytuacja demograficzna Polski. Raport 2023-2024”,<br>Główny Urząd Statystyczny</p<p>Mieszkanie w małej miejscowości wiąże się z szeregiem problemów wywołanych na przykład brakiem dostępu do edukacji wyższej czy możliwości pracy w wymarzonym miejscu. W sytuacjach, gdy rodzina jest wielodzietna, dzieci często muszą też ustalać między sobą, kto zostanie z rodzicami, żeby opiekować się nimi na starość lub przejąć ich obowiązki zawodowe. Moja rozmówczyni, Ola, mówi o swoim bracie tak: ,,No ja wiem, że on się nie wyprowadzi, on po prostu tam zostanie i będzie się opiekował rodzicami.”</p><p>W rozmowie z Olą starałam się dowiedzieć, jak wygląda perspektywa młodej dorosłej, która mieszka w Mierzeszynie, czyli we wsi oddalonej o 30 km od Gdańska. Jej życie skupia się głównie wokół  Trójmiasta. Tam ukończyła wyższą edukację, pracuje i planuje większość swoich aktywności. Rozmawiałyśmy o problemach, z jakimi boryka się osoba w jej sytuacji oraz możliwych utrudnieniach dla innych mieszkańców wsi.</p>'
$ws.Range("E3").WrapText = $false

$ws.Range("E11").Value = '<b>Bazarek Vinted</b><br><br><p>Okolice pętli Metro Wilanowska stały się obszarem moich obserwacji i zainteresowania, od kiedy codziennie przesiadam się tam z autobusu na metro i z powrotem. Gdy pogoda sprzyja kilka starszych pań wystawia maleńkie bazarki, na którym sprzedają ubrania, książki, rośliny i inne drobiazgi. Kiedyś kupiłam od jednej pani roślinę doniczkową, do dzisiaj zdobi mój pokój. Taka forma sprzedaży kojarzy mi się z internetową platformą Vinted, na której prywatne osoby, mogą sprzedawać używaną odzież i akcesoria do domu. Dlatego ten powstający na kilka godzin bazarek, nazywam bazarkiem Vinted.</p><p> Ciekawią mnie historie osób sprzedających, jak wygląda ich dzień, czemu sprzedają... Wstępnie miałam okazję porozmawiać z jedną panią, niestety finalnie nie zgodziła się udzielić wywiadu na wystawę. A od kiedy przyszły mrozy, panie nie rozstawiają straganów, więc nie miałam okazji się zaprosić innej pani.</p><p>Tym sposobem moim rozmówcą został codzienny przechodzień w tym miejscu. Tak jak ja przesiada się z jednego środku transportu na drugi. Od dziewięciu lat mija bazarki, czasem się zatrzymuje, obserwuje. Dzieli się swoimi spostrzeżeniami.</p>'
$ws.Range("E11").WrapText = $false

$ws.Range("E13").Value = '<b>Gorzko, gorzko!</b><br><br><p>Śluby i wesela są wydarzeniami, które szczególnie mocno zapisują się w pamięci lokalnej społeczności, dlatego stały się dla mnie punktem wyjścia do opowieści o miejscu, relacjach i zmianach społecznych. Rozmowy z moją prababcią i mamą pozwoliły mi zobaczyć, że lokalność to coś żywego, zmieniającego się, ale i zakorzenionego w tych samych przestrzeniach. Zależało mi na poznaniu emocji i wspomnień, aby poprzez indywidualne historie wokół ślubów i wesel opowiedzieć o wspólnym doświadczeniu mieszkańców gminy Wolanów. Wybrałam ją jako miejsce mi bliskie, znane z osobistego doświadczenia i rodzinnych opowieści, co pozwoliło mi spojrzeć na lokalność z perspektywy wewnętrznej i międzypokoleniowej.</p>'
$ws.Range("E13").WrapText = $false

$ws.Range("E14").Value = '<b>Powrót do lokalności</b><br><br><p>Rozmowę z Panem Robertem przeprowadziłam z potrzeby osobistej ciekawości i nostalgii związanej z miejscowością Przytyk. Chciałam dowiedzieć się więcej o miejscu, które pomimo dziecięcych wspomnień, wydaje mi się zupełnie obce - zrozumieć jak wyglądało ono kiedyś oraz jak dziś postrzegane jest przez tych, którzy do niego wracają.</p><p>Pan Robert, 52-letni mieszkaniec Przytyka, spędził tutaj dzieciństwo i młodość, a następnie przez trzy dekady żył poza jego granicami - w innych obszarach Polski, ale również w różnych krajach Europy. Jego powrót nie jest sentymentalną podróżą „do korzeni”, ale konfrontacją z miejscem, które przez lata przestało być takie, jakim je zapamiętał. W jego słowach pojawia się zarówno wspomnienie dziecięcej wolności, jak i dystans wobec współczesnej rzeczywistości małej miejscowości.</p><p>Wybór tej historii pozwolił mi przyjrzeć się, czym dziś może być „lokalność” - nie jako przestrzeń geograficzna, ale ta emocjonalna i biograficzna. Choć Pan Robert wrócił, w jego opowieści pojawia się pragnienie mobilności, marzenie o życiu w kamperze, o wolności od przywiązań i rzeczy. To opowieść o przemijaniu, powrocie z konieczności i o próbie odnalezienia siebie w miejscu, które tak jak wszystko inne, zmieniło się.</p>'
$ws.Range("E14").WrapText = $false

# Remove the stray "iujhyrwe" placeholder value that was left in column E (row 12).
$ws.Range("E12").ClearContents()

# Update the remembered selection/active cell for the sheet.
$ws.Range("E18").Select()
